# Add a new "2022-Q1" sheet (before the "总计" sheet) with fund holding
# data, and update the "总计" (totals) summary sheet with a new row for
# 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, positioned just before "总计".
#    NOTE: Worksheets.Add() renumbers sheet positions, which invalidates
#    any previously-captured sheet handle. So do this FIRST, before
#    grabbing any other sheet references.
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)   # "总计" - currently last
$q1Sheet = $wb.Worksheets.Add($lastSheet)
$q1Sheet.Name = "2022-Q1"

# Now that no more sheet-add/remove operations remain, it is safe to
# fetch (and keep) the other sheet handles we need.
$q4Sheet = $wb.Worksheets.Item("2021-Q4")       # style donor
$totalSheet = $wb.Worksheets.Item("总计")

# Copy header/number formatting (borders, bold, centering) from the
# "2021-Q4" sheet, which already carries the desired style.
$q4Sheet.Range("B1:H1").Copy()
$q1Sheet.Range("B1:H1").PasteSpecial(-4122)
$q4Sheet.Range("A2").Copy()
$q1Sheet.Range("A2:A7").PasteSpecial(-4122)

# Header row
$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"

# Data rows. Column A / H are numeric; B-G are text (leading zeros /
# trailing zeros in the source data must be preserved verbatim), so
# they are written with a leading apostrophe to force text storage.
$q1Data = @(
    @(0, "002910", "易方达供给改革灵活配置混合", "49.29", "87.54", "3.56", "1.7547", 10),
    @(1, "002281", "建信裕利灵活配置混合",       "1.10",  "88.94", "6.19", "0.0681", 1),
    @(2, "009124", "华泰保兴科荣混合A",          "5.12",  "22.53", "1.32", "0.0676", 7),
    @(3, "002378", "建信弘利灵活配置混合",       "1.03",  "89.57", "5.97", "0.0615", 2),
    @(4, "007385", "华泰保兴安盈三个月定期开放混合", "7.13", "21.22", "0.63", "0.0449", 8),
    @(5, "009125", "华泰保兴科荣混合C",          "0.00",  "22.53", "1.32", 0, 7)
)

$r = 2
foreach ($row in $q1Data) {
    $q1Sheet.Cells.Item($r, 1).Value = $row[0]
    $q1Sheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $q1Sheet.Cells.Item($r, 3).Value = "'" + $row[2]
    $q1Sheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $q1Sheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $q1Sheet.Cells.Item($r, 6).Value = "'" + $row[5]
    if ($r -eq 7) {
        $q1Sheet.Cells.Item($r, 7).Value = 0
    } else {
        $q1Sheet.Cells.Item($r, 7).Value = "'" + $row[6]
    }
    $q1Sheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Add a new top data row to "总计" for 2022-Q1, shifting the
#    existing 2021-Q4 / 2021-Q3 rows down.
# ---------------------------------------------------------------------
$totalSheet.Rows(2).Insert()

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 6
$totalSheet.Cells.Item(2, 4).Value = 2

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2

Write-Output "done"
